$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "Your company is developing an app to digitally restore old photographs. You notice that many of these images contain salt-and-pepper noise. You want to apply noise reduction while preserving the image's edges.What should you do?",
        "ques_type": 2,
        "options": [
            "Use a median filter.",
            "Use a Gaussian filter.",
            "Apply histogram equalization.",
            "Increase the image contrast."
        ],
        "score": "Use a median filter."
    },
    {
        "title": "You have deployed an object detection system in a crowded city street. However, the system often detects multiple bounding boxes for the same object, leading to redundant detections. What should you do?",
        "ques_type": 2,
        "options": [
            "Apply non-maximum suppression (NMS).",
            "Increase the threshold for detection.",
            "Reduce the number of anchor boxes.",
            "Apply Gaussian blur to the input image."
        ],
        "score": "Apply non-maximum suppression (NMS)."
    },
    {
        "title": "You are developing a system to identify tumors in medical images. The system needs to classify each pixel in the image to either 'tumor' or 'non-tumor.'Which architecture should you use?",
        "ques_type": 2,
        "options": [
            "U-Net",
            "YOLO",
            "RCNN",
            "SSD"
        ],
        "score": "U-Net"
    },
    {
        "title": "Your company is building a tool for real-time emotion recognition using facial features. You want to capture subtle variations in facial features across different regions of the face.What should you do to effectively capture these regional variations?",
        "ques_type": 2,
        "options": [
            "Use local binary patterns (LBPs).",
            "Apply global histogram equalization.",
            "Use edge detection methods.",
            "Apply Fourier transform."
        ],
        "score": "Use local binary patterns (LBPs)."
    }
]
'@

# Clear existing contents/formatting of A1 and A2 (A1 currently holds a bordered/bold "0",
# A2 currently holds the long question text). We collapse everything down to a single
# plain A1 cell containing the reformatted text.
$ws.Cells.Clear()

$ws.Range("A1").Value = $newText

# The new text contains embedded line breaks, which auto-expands the row height.
# Re-fit the row so it doesn't keep an explicit custom height.
$ws.Rows(1).AutoFit()
